# Update of Rachel setup
# Applies the changes from the commit: adds a new row 53 (sequential
# numbers 1-50 across columns B:AY, mirroring the header row 2), and
# updates the active sheet view (zoom + selected cell) to match the
# author's last-saved window state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 53: plain sequential integers 1..50 in B53:AY53 ---------------
for ($i = 0; $i -lt 50; $i++) {
    $ws.Cells.Item(53, 2 + $i).Value = $i + 1
}

# --- Sheet view / window state ---------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 120
$null = $ws.Range("C48").Select()

Write-Host "Row 53 populated and view updated"
